$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fill in the new "Will Do" / "Next activities" cells for the ThanhLC row (row 4)
$ws.Range("B4").Value = "Watch tutorials on CSS3, HTML5"
$ws.Range("C4").Value = "Continue on watching tutorials"

# Update the saved selection to C3
$ws.Range("C3").Select()

# Update the absolute path recorded for the workbook (Microsoft 365 x15ac:absPath)
$wb.Path = "C:\Users\mGlushed\Documents\_GREENWICH\Web\GitHub\Repo\Meeting Minutes\"
